$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Add the new version row (row 4) with the FSAddress changelog entry
$ws.Range("A4").Value = "[1.5]"
$ws.Range("B4").Value = "add FSAddress field to the templates and csv and code"
$ws.Range("C4").Value = 43222
$ws.Range("C4").NumberFormat = "d-mmm-yy"

# Grow Table2 to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C5"))

# Move the selection the way the author left it
$ws.Range("A5").Select() | Out-Null
